# Update "presentase fuzzy matching.xlsx" - revise one of the fuzzy
# matching accuracy scores (B10) in the first results table; the
# dependent average (B14) and percentage-change formula (F3) recalc
# automatically. Also move the active selection to B19, matching the
# saved cursor position in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = 0.9397

$ws.Range("B19").Select()
